$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row data for rows 9-59 (inclusive): each entry lists the A,B,C,D,E,F values for that row
$rowData = @{
    9 = @(7, 0, 0, 31, 10, 21)
    10 = @(8, 0, 0, 38, 17, 21)
    11 = @(9, 0, 0, 43, 22, 21)
    12 = @(10, 0, 0, 50, 29, 21)
    13 = @(11, 0, 21, 59, 38, 0)
    14 = @(12, 0, 21, 61, 34, 6)
    15 = @(13, 0, 21, 68, 32, 15)
    16 = @(14, 0, 21, 77, 35, 21)
    17 = @(15, 0, 21, 82, 35, 26)
    18 = @(16, 0, 21, 89, 32, 36)
    19 = @(17, 0, 27, 106, 45, 34)
    20 = @(18, 0, 36, 121, 55, 30)
    21 = @(19, 0, 42, 144, 68, 34)
    22 = @(20, 0, 47, 157, 76, 34)
    23 = @(21, 2, 57, 183, 94, 32)
    24 = @(22, 2, 61, 202, 103, 38)
    25 = @(23, 2, 66, 243, 127, 50)
    26 = @(24, 2, 76, 296, 162, 58)
    27 = @(25, 2, 81, 356, 204, 71)
    28 = @(26, 2, 89, 409, 232, 88)
    29 = @(27, 3, 99, 478, 284, 95)
    30 = @(28, 4, 116, 547, 332, 99)
    31 = @(29, 4, 134, 627, 353, 140)
    32 = @(30, 4, 152, 728, 399, 177)
    33 = @(31, 7, 177, 844, 461, 206)
    34 = @(32, 7, 194, 971, 527, 250)
    35 = @(33, 7, 215, 1122, 619, 288)
    36 = @(34, 7, 274, 1276, 686, 316)
    37 = @(35, 8, 329, 1479, 805, 345)
    38 = @(36, 8, 383, 1695, 898, 414)
    39 = @(37, 9, 444, 1923, 1004, 475)
    40 = @(38, 11, 503, 2167, 1116, 548)
    41 = @(39, 12, 590, 2472, 1282, 600)
    42 = @(40, 14, 674, 2806, 1431, 701)
    43 = @(41, 17, 797, 3220, 1632, 791)
    44 = @(42, 18, 919, 3704, 1886, 899)
    45 = @(43, 19, 1051, 4187, 2137, 999)
    46 = @(44, 25, 1190, 4732, 2418, 1124)
    47 = @(45, 27, 1375, 5368, 2722, 1271)
    48 = @(46, 32, 1588, 6103, 3075, 1440)
    49 = @(47, 38, 1818, 6926, 3475, 1633)
    50 = @(48, 42, 2050, 7846, 3888, 1908)
    51 = @(49, 53, 2314, 8881, 4438, 2129)
    52 = @(50, 62, 2646, 9996, 4948, 2402)
    53 = @(51, 73, 3028, 11269, 5507, 2734)
    54 = @(52, 82, 3451, 12673, 6153, 3069)
    55 = @(53, 91, 3958, 14147, 6740, 3449)
    56 = @(54, 109, 4443, 15775, 7408, 3924)
    57 = @(55, 120, 5048, 17574, 8109, 4417)
    58 = @(56, 143, 5762, 19335, 8689, 4884)
    59 = @(57, 164, 6520, 21333, 9351, 5462)
}

foreach ($r in ($rowData.Keys | Sort-Object {[int]$_})) {
    $rowNum = [int]$r
    $vals = $rowData[$r]
    for ($col = 1; $col -le 6; $col++) {
        $ws.Cells.Item($rowNum, $col).Value = $vals[$col - 1]
    }
}

# Row 59 is brand new. Give column A the same style (bold / bordered / centered)
# used by the other day-number cells in column A (e.g. A58).
$ws.Cells.Item(58, 1).Copy() | Out-Null
$ws.Cells.Item(59, 1).PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = 0

# Re-assert the value for A59 in case the paste-formats touched it.
$ws.Cells.Item(59, 1).Value = 57
